$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lembar1")

# Update values on row 2
$ws.Range("A2").Value = "test"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "Bogor"

# Clear out rows 3-7 (columns B:D), keep column A style but blank value
$ws.Range("A3:D7").ClearContents() | Out-Null

# Reset the sheet view: remove the frozen/scrolled topLeftCell and move selection to D2
$ws.Activate() | Out-Null
$ws.Range("D2").Select() | Out-Null
